$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 542.2941
$ws.Range("J33").Value2 = 1150
$ws.Range("L33").Value2 = 1150
$ws.Range("N33").Value2 = -1608
$ws.Range("H43").Value2 = 2416.3333
$ws.Range("I43").Value2 = 2099.6
$ws.Range("J43").Value2 = 4000
$ws.Range("K43").Value2 = 2099.6
$ws.Range("L43").Value2 = 4000
$ws.Range("M43").Value2 = -2030.6
$ws.Range("N43").Value2 = -4138
$ws.Range("H64").Value2 = 14605.066
$ws.Range("J64").Value2 = 9047.700000000001
$ws.Range("L64").Value2 = 9047.700000000001
$ws.Range("N64").Value2 = -9543.700000000001
$ws.Range("H67").Value2 = 14605.066
$ws.Range("J67").Value2 = 9047.700000000001
$ws.Range("L67").Value2 = 9047.700000000001
$ws.Range("N67").Value2 = -10763.7
$ws.Range("H86").Value2 = 5779.143
$ws.Range("I86").Value2 = 9500
$ws.Range("J86").Value2 = 4290.8
$ws.Range("K86").Value2 = 9500
$ws.Range("L86").Value2 = 4290.8
$ws.Range("M86").Value2 = -8377
$ws.Range("N86").Value2 = -6536.8
$ws.Range("H89").Value2 = 5779.143
$ws.Range("I89").Value2 = 9500
$ws.Range("J89").Value2 = 4290.8
$ws.Range("K89").Value2 = 47500
$ws.Range("L89").Value2 = 21454
$ws.Range("M89").Value2 = -41884
$ws.Range("N89").Value2 = -32686
$ws.Range("H111").Value2 = 6024.2144
$ws.Range("I111").Value2 = 4756.5
$ws.Range("K111").Value2 = 14269.5
$ws.Range("M111").Value2 = -11202.5
$ws.Range("H123").Value2 = 59998.168
$ws.Range("J123").Value2 = 59998.168
$ws.Range("L123").Value2 = 59998.168
$ws.Range("N123").Value2 = -69798.16800000001
$ws.Range("H137").Value2 = 1355.5476
$ws.Range("I137").Value2 = 897.8
$ws.Range("J137").Value2 = 2499.9167
$ws.Range("K137").Value2 = 2693.4
$ws.Range("L137").Value2 = 7499.750100000001
$ws.Range("M137").Value2 = -143.3999999999996
$ws.Range("N137").Value2 = -12599.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value2 = 0
$ws.Range("J24").Value2 = 0
$ws.Range("H61").Value2 = 3357.25
$ws.Range("I61").Value2 = 1649.6428
$ws.Range("K61").Value2 = 1649.6428
$ws.Range("M61").Value2 = -1437.6428
$ws.Range("H74").Value2 = 1467.4375
$ws.Range("I74").Value2 = 1326
$ws.Range("J74").Value2 = 1778.6
$ws.Range("K74").Value2 = 1326
$ws.Range("L74").Value2 = 1778.6
$ws.Range("M74").Value2 = -452
$ws.Range("N74").Value2 = -3526.6
$ws.Range("H77").Value2 = 1467.4375
$ws.Range("I77").Value2 = 1326
$ws.Range("J77").Value2 = 1778.6
$ws.Range("K77").Value2 = 6630
$ws.Range("L77").Value2 = 8893
$ws.Range("M77").Value2 = -2262
$ws.Range("N77").Value2 = -17629
$ws.Range("H100").Value2 = 0
$ws.Range("J100").Value2 = 0
$ws.Range("H101").Value2 = 38000
$ws.Range("J101").Value2 = 38000
$ws.Range("L101").Value2 = 38000
$ws.Range("N101").Value2 = -44490
$ws.Range("H110").Value2 = 2840.111
$ws.Range("I110").Value2 = 2840.111
$ws.Range("K110").Value2 = 2840.111
$ws.Range("M110").Value2 = -795.1109999999999
$ws.Range("H132").Value2 = 3698.2693
$ws.Range("I132").Value2 = 3225.7856
$ws.Range("K132").Value2 = 9677.356800000001
$ws.Range("M132").Value2 = -7147.356800000001
$ws.Range("H135").Value2 = 55783.766
$ws.Range("J135").Value2 = 55783.766
$ws.Range("L135").Value2 = 55783.766
$ws.Range("N135").Value2 = -65923.766
$ws.Range("H136").Value2 = 3357.25
$ws.Range("I136").Value2 = 1649.6428
$ws.Range("K136").Value2 = 4948.928400000001
$ws.Range("M136").Value2 = -2398.928400000001
$ws.Range("H139").Value2 = 79000
$ws.Range("J139").Value2 = 78500
$ws.Range("L139").Value2 = 78500
$ws.Range("N139").Value2 = -88780
$ws.Range("N24").ClearContents()
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 4817.4
$ws.Range("I20").Value2 = 2779.5
$ws.Range("K20").Value2 = 2779.5
$ws.Range("M20").Value2 = -2532.5
$ws.Range("H134").Value2 = 2644.8367
$ws.Range("I134").Value2 = 1776.119
$ws.Range("J134").Value2 = 7857.143
$ws.Range("K134").Value2 = 5328.357
$ws.Range("L134").Value2 = 23571.429
$ws.Range("M134").Value2 = -2793.357
$ws.Range("N134").Value2 = -28641.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 628.8
$ws.Range("I16").Value2 = 471.42856
$ws.Range("J16").Value2 = 996
$ws.Range("K16").Value2 = 471.42856
$ws.Range("L16").Value2 = 996
$ws.Range("M16").Value2 = -184.42856
$ws.Range("N16").Value2 = -1570
$ws.Range("H22").Value2 = 500
$ws.Range("I22").Value2 = 0
$ws.Range("J22").Value2 = 500
$ws.Range("K22").Value2 = 0
$ws.Range("N22").Value2 = -1200
$ws.Range("H74").Value2 = 62430.6
$ws.Range("J74").Value2 = 62430.6
$ws.Range("L74").Value2 = 62430.6
$ws.Range("N74").Value2 = -64178.6
$ws.Range("H77").Value2 = 62430.6
$ws.Range("J77").Value2 = 62430.6
$ws.Range("L77").Value2 = 187291.8
$ws.Range("N77").Value2 = -196027.8
$ws.Range("H105").Value2 = 31260750
$ws.Range("I105").Value2 = 1833
$ws.Range("K105").Value2 = 1833
$ws.Range("M105").Value2 = -86
$ws.Range("H113").Value2 = 628.8
$ws.Range("I113").Value2 = 471.42856
$ws.Range("J113").Value2 = 996
$ws.Range("K113").Value2 = 471.42856
$ws.Range("L113").Value2 = 996
$ws.Range("M113").Value2 = 1698.57144
$ws.Range("N113").Value2 = -5336
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value2 = 1800
$ws.Range("I17").Value2 = 1400
$ws.Range("J17").Value2 = 2000
$ws.Range("K17").Value2 = 4200
$ws.Range("L17").Value2 = 6000
$ws.Range("M17").Value2 = -4031
$ws.Range("N17").Value2 = -6338
$ws.Range("H60").Value2 = 700.8
$ws.Range("I60").Value2 = 310
$ws.Range("K60").Value2 = 930
$ws.Range("M60").Value2 = -679
$ws.Range("H70").Value2 = 150583
$ws.Range("I70").Value2 = 1150
$ws.Range("K70").Value2 = 3450
$ws.Range("M70").Value2 = -3135
$ws.Range("H73").Value2 = 150583
$ws.Range("I73").Value2 = 1150
$ws.Range("K73").Value2 = 3450
$ws.Range("M73").Value2 = -2358
$ws.Range("H139").Value2 = 2938.7778
$ws.Range("I139").Value2 = 2938.7778
$ws.Range("K139").Value2 = 8816.3334
$ws.Range("M139").Value2 = -3676.3334
$ws.Range("H140").Value2 = 1780.9286
$ws.Range("I140").Value2 = 1687.1538
$ws.Range("K140").Value2 = 5061.4614
$ws.Range("M140").Value2 = 118.5385999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 328.12
$ws.Range("I2").Value2 = 434.8
$ws.Range("K2").Value2 = 434.8
$ws.Range("M2").Value2 = -321.8
$ws.Range("H33").Value2 = 7903.6665
$ws.Range("J33").Value2 = 7903.6665
$ws.Range("L33").Value2 = 7903.6665
$ws.Range("N33").Value2 = -8407.666499999999
$ws.Range("H107").Value2 = 541.3
$ws.Range("I107").Value2 = 310.25
$ws.Range("J107").Value2 = 695.3333
$ws.Range("K107").Value2 = 310.25
$ws.Range("L107").Value2 = 695.3333
$ws.Range("M107").Value2 = 1609.75
$ws.Range("N107").Value2 = -4535.3333
$ws.Range("H113").Value2 = 11071
$ws.Range("J113").Value2 = 13947.5
$ws.Range("L113").Value2 = 13947.5
$ws.Range("N113").Value2 = -18287.5
$ws.Range("H132").Value2 = 3508.4348
$ws.Range("I132").Value2 = 1959.3572
$ws.Range("J132").Value2 = 5918.1113
$ws.Range("K132").Value2 = 5878.071599999999
$ws.Range("L132").Value2 = 17754.3339
$ws.Range("M132").Value2 = -3348.071599999999
$ws.Range("N132").Value2 = -22814.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 8105.8125
$ws.Range("I81").Value2 = 10266.083
$ws.Range("J81").Value2 = 1625
$ws.Range("K81").Value2 = 20532.166
$ws.Range("L81").Value2 = 3250
$ws.Range("M81").Value2 = -19471.166
$ws.Range("N81").Value2 = -5372
$ws.Range("H84").Value2 = 8105.8125
$ws.Range("I84").Value2 = 10266.083
$ws.Range("J84").Value2 = 1625
$ws.Range("K84").Value2 = 102660.83
$ws.Range("L84").Value2 = 16250
$ws.Range("M84").Value2 = -97356.83
$ws.Range("N84").Value2 = -26858
$ws.Range("H107").Value2 = 354.7857
$ws.Range("I107").Value2 = 372.25
$ws.Range("K107").Value2 = 1116.75
$ws.Range("M107").Value2 = 803.25
$ws.Range("H122").Value2 = 3243.3784
$ws.Range("I122").Value2 = 2845.84
$ws.Range("J122").Value2 = 4071.5833
$ws.Range("K122").Value2 = 8537.52
$ws.Range("L122").Value2 = 12214.7499
$ws.Range("M122").Value2 = -6087.52
$ws.Range("N122").Value2 = -17114.7499
$ws.Range("H132").Value2 = 1391.3438
$ws.Range("I132").Value2 = 933.87036
$ws.Range("K132").Value2 = 2801.61108
$ws.Range("M132").Value2 = -271.6110800000001
$ws.Range("H136").Value2 = 1620.8448
$ws.Range("I136").Value2 = 832.73914
$ws.Range("J136").Value2 = 4641.9165
$ws.Range("K136").Value2 = 2498.21742
$ws.Range("L136").Value2 = 13925.7495
$ws.Range("M136").Value2 = 51.78258000000005
$ws.Range("N136").Value2 = -19025.7495
